# Append a new row (row 4) of data to Sheet1, mirroring the existing rows'
# inline-text formatting (MIGRATION DATE, FINANCIAL INSTITUTION NAME,
# ENTITY ID, ADDRESS).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-looking string ("2025-10-17") that must be stored
# as literal text, not auto-converted to a date serial number. Pre-format
# the cell as Text before writing, then clear the formatting override so
# the cell ends up with no explicit style (matching the unstyled cells in
# rows 2 and 3).
$dateCell = $ws.Range("A4")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-10-17"
$dateCell.ClearFormats()

$ws.Range("B4").Value = "ZZZ"
$ws.Range("C4").Value = "456CDX009"
$ws.Range("D4").Value = "Anna Nagar"
